$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New formula cells
$ws.Range("C8").Formula = "=220/344"
$ws.Range("C9").Formula = "=220/(99+220)"

# Column F width
$ws.Columns.Item(6).ColumnWidth = 15.86

# Legend table
$ws.Range("E21").Value = "P"
$ws.Range("F21").Value = "N"
$ws.Range("D22").Value = "T"
$ws.Range("E22").Value = "Correct answers "
$ws.Range("F22").Value = "No wikpedia page"
$ws.Range("D23").Value = "F"
$ws.Range("E23").Value = "Wrong answers "
$ws.Range("F23").Value = "Missing link"

# Styling: D21 white fill (no font change -- inherits plain font)
$ws.Range("B7").Copy()
$ws.Range("D21").PasteSpecial(-4122)
$ws.Range("D21").ClearContents()
$ws.Range("D21").Interior.Color = 16777215

# E21 green fill, bold white font  (first bold font -> font index 2)
$ws.Range("E21").Interior.Color = 65280
$ws.Range("E21").Font.Bold = $true
$ws.Range("E21").Font.Color = 16777215

# F21 red fill, bold light-gray font  (second bold font -> font index 3)
$ws.Range("F21").Interior.Color = 255
$ws.Range("F21").Font.Bold = $true
$ws.Range("F21").Font.Color = 15987699

# D22 green fill, bold black font  (third bold font -> font index 4)
$ws.Range("D22").Interior.Color = 65280
$ws.Range("D22").Font.Bold = $true
$ws.Range("D22").Font.Color = 0

# E22,F22 white fill bold black font
$ws.Range("E22:F22").Interior.Color = 16777215
$ws.Range("E22:F22").Font.Bold = $true
$ws.Range("E22:F22").Font.Color = 0

# D23 red fill bold black font
$ws.Range("D23").Interior.Color = 255
$ws.Range("D23").Font.Bold = $true
$ws.Range("D23").Font.Color = 0

# E23,F23 white fill bold black font
$ws.Range("E23:F23").Interior.Color = 16777215
$ws.Range("E23:F23").Font.Bold = $true
$ws.Range("E23:F23").Font.Color = 0

# Move/resize the chart
$chart = $ws.ChartObjects(1)
$chart.Chart.ChartArea.Left = 0
